# Tim Seifert.xlsx - turn the generic "Sheet1" into a per-match scrape row
# that also records which match number the innings belongs to.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "Tim Seifert"
$ws.Name = "Tim Seifert"

# Insert a brand-new column A ("matchNo"), shifting teamName..result one
# column to the right (B..M instead of A..L).
$ws.Columns.Item(1).Insert(-4161)  # -4161 = xlShiftToRight

# Populate the new column's header + the single data row's value.
$ws.Range("A1").Value = "matchNo"
$ws.Range("A2").Value = "45th"
